# Auto-generated Excel COM-interop script
# Applies cell value updates across multiple sheets per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H74" = 10222.111
    "I74" = 7999.6665
    "K74" = 7999.6665
    "M74" = -7063.6665
    "H77" = 10222.111
    "I77" = 7999.6665
    "K77" = 43165
    "M77" = -35318.3325
    "H100" = 4189.8
    "I100" = 3649.6667
    "K100" = 3649.6667
    "M100" = -3108.6667
    "H112" = 3513.818
    "J112" = 3666.2
    "L112" = 10998.6
    "N112" = -13214.6
    "H137" = 2566.2856
    "I137" = 2410.8333
    "K137" = 7232.499899999999
    "M137" = -4682.499899999999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H32" = 4099.3716
    "I32" = 4099.3716
    "K32" = 4099.3716
    "M32" = -3812.3716
    "H37" = 16332.667
    "I37" = 4999
    "J37" = 21999.5
    "K37" = 4999
    "L37" = 21999.5
    "M37" = -4726
    "N37" = -22545.5
    "H57" = 50000
    "I57" = 50000
    "K57" = 50000
    "M57" = -49246
    "H61" = 2762.1765
    "I61" = 1351.1111
    "K61" = 1351.1111
    "M61" = -1139.1111
    "H102" = 2224.2
    "I102" = 2113.3333
    "J102" = 2390.5
    "K102" = 2113.3333
    "L102" = 2390.5
    "M102" = -491.3332999999998
    "N102" = -5634.5
    "H136" = 2762.1765
    "I136" = 1351.1111
    "K136" = 4053.3333
    "M136" = -1503.3333
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H20" = 4449.625
    "I20" = 4228.143
    "K20" = 4228.143
    "M20" = -3981.143
    "H105" = 2849.5
    "I105" = 2819.4
    "J105" = 3000
    "K105" = 2819.4
    "L105" = 3000
    "M105" = -1072.4
    "N105" = -6494
    "H135" = 40000
    "I135" = 40000
    "K135" = 40000
    "M135" = -34930
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H7" = 367
    "I7" = 139.4
    "K7" = 139.4
    "M7" = -26.40000000000001
    "H62" = 3622
    "I62" = 3424.2856
    "K62" = 3424.2856
    "M62" = -2800.2856
    "H65" = 3622
    "I65" = 3424.2856
    "K65" = 17121.428
    "M65" = -14001.428
    "H86" = 12003.917
    "I86" = 6262.3335
    "J86" = 13917.777
    "K86" = 6262.3335
    "L86" = 13917.777
    "M86" = -5139.3335
    "N86" = -16163.777
    "H89" = 12003.917
    "I89" = 6262.3335
    "J89" = 13917.777
    "K89" = 31311.6675
    "L89" = 69588.88499999999
    "M89" = -25695.6675
    "N89" = -80820.88499999999
    "H94" = 1692.3334
    "I94" = 1292
    "J94" = 1892.5
    "K94" = 1292
    "L94" = 1892.5
    "M94" = -841
    "N94" = -2794.5
    "H132" = 3155.25
    "I132" = 2049
    "K132" = 6147
    "M132" = -3617
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H68" = 1800
    "I68" = 0
    "K68" = 0
    "H71" = 1800
    "I71" = 0
    "K71" = 0
    "H80" = 0
    "I80" = 0
    "K80" = 0
    "H83" = 0
    "I83" = 0
    "K83" = 0
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$deletes = @("M68", "M71", "M80", "M83")
foreach ($cellRef in $deletes) {
    $ws.Range($cellRef).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H70" = 6999
    "J70" = 6998
    "L70" = 6998
    "N70" = -7538
    "H73" = 6999
    "J73" = 6998
    "L73" = 6998
    "N73" = -8870
    "H102" = 2366.6667
    "I102" = 2366.6667
    "K102" = 2366.6667
    "M102" = -744.6667000000002
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H22" = 7499
    "I22" = 1998.3334
    "J22" = 9856.429
    "K22" = 1998.3334
    "L22" = 9856.429
    "M22" = -1703.3334
    "N22" = -10446.429
    "H27" = 7499
    "I27" = 1998.3334
    "J27" = 9856.429
    "K27" = 1998.3334
    "L27" = 9856.429
    "M27" = -1891.3334
    "N27" = -10070.429
    "H46" = 2280.5454
    "I46" = 1583.7142
    "K46" = 1583.7142
    "M46" = -1395.7142
    "H93" = 1641.8
    "I93" = 1641.8
    "K93" = 1641.8
    "M93" = -393.8
    "H100" = 1318
    "I100" = 1318
    "J100" = 0
    "K100" = 1318
    "L100" = 0
    "M100" = -777
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$deletes = @("N100")
foreach ($cellRef in $deletes) {
    $ws.Range($cellRef).ClearContents()
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$updates = @{
    "H81" = 20581.666
    "I81" = 22654.375
    "J81" = 4000
    "K81" = 45308.75
    "L81" = 8000
    "M81" = -44247.75
    "N81" = -10122
    "H84" = 20581.666
    "I84" = 22654.375
    "J84" = 4000
    "K84" = 226543.75
    "L84" = 40000
    "M84" = -221239.75
    "N84" = -50608
    "H132" = 2120.1304
    "I132" = 1931.4117
    "J132" = 2654.8333
    "K132" = 5794.2351
    "L132" = 7964.499899999999
    "M132" = -3264.2351
    "N132" = -13024.4999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Host "Applied all Sephirot_Profits updates"